$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Source" block at the bottom of the sheet (rows 47-56) is being
# restructured:
#   - a new blank line is inserted right after "Source:" (row 47)
#   - the old SERCOTEC URL cell (which carried a live hyperlink) is
#     replaced by a plain-text cell one row further down, and the
#     hyperlink itself is removed
#   - the two long citation paragraphs are shortened to just "SERCOTEC"

# 1) Insert a new blank row right after row 47 ("Source:").
#    This pushes the old rows 48-56 down to 49-57.
$ws.Rows(48).Insert()

# 2) The old hyperlinked URL row (now at row 50) is removed...
$ws.Rows(50).Delete()

# 3) ...and re-created two rows further down as a plain text cell
#    (no hyperlink), using the same "source" style as its neighbours.
$ws.Rows(51).Insert()
$ws.Range("A51").Value = "http://www.sercotec.cl/Portals/0/MANUALES/situaci%C3%B3n%20de%20la%20microempresa.pdf"
$ws.Range("A51").Style = "source"

# 4) Shorten the two long citation paragraphs down to "SERCOTEC".
$ws.Range("A55").Value = "SERCOTEC"
$ws.Range("A57").Value = "SERCOTEC"

# 5) The sheet no longer has any live hyperlinks.
$ws.Hyperlinks.Delete()
